$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "Daily Orders": a brand-new incoming order (Pooja / Wheat Chapati)
# arrives and is inserted as the new row 2, pushing the existing
# "Anuradha N" order down to row 3.
# -------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daily Orders")

$ws1.Rows.Item(2).Insert()

$ws1.Range("A2").Value = 2
$ws1.Range("B2").Value = "2026-01-13 10:20"
$ws1.Range("C2").Value = "Pooja"
$ws1.Range("D2").Value = "A 1608"

# Phone column is blank for this order - format as text first so the
# empty value still materialises a (blank) cell instead of being dropped.
$ws1.Range("E2").NumberFormat = "@"
$ws1.Range("E2").Value = ""

$ws1.Range("F2").Value = "Wheat Chapati x1"
$ws1.Range("G2").Value = 15
$ws1.Range("H2").Value = "NEW"
$ws1.Range("I2").Value = "PENDING"

# Collection Date is a plain "yyyy-mm-dd" looking string that must stay
# text, not get auto-converted into a real date serial number.
$ws1.Range("J2").NumberFormat = "@"
$ws1.Range("J2").Value = "2026-01-13"

$ws1.Range("K2").Value = "15:50"

# Notes / Cancel Reason / Feedback are blank for a brand new order.
$ws1.Range("L2").NumberFormat = "@"
$ws1.Range("L2").Value = ""
$ws1.Range("M2").NumberFormat = "@"
$ws1.Range("M2").Value = ""
$ws1.Range("N2").NumberFormat = "@"
$ws1.Range("N2").Value = ""

# -------------------------------------------------------------------------
# Sheet "Summary": roll the new order into the daily totals.
# -------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("A2").Value = 2   # Total Orders: 1 -> 2
$ws2.Range("B2").Value = 1   # New: 0 -> 1
$ws2.Range("G2").Value = 45  # Total Revenue: 30 -> 45

# -------------------------------------------------------------------------
# Sheet "Items Breakdown": add the new item line above the existing one.
# -------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Items Breakdown")

$ws3.Rows.Item(2).Insert()

$ws3.Range("A2").Value = "Wheat Chapati"
$ws3.Range("B2").Value = 1
$ws3.Range("C2").Value = 15
